# Update "Name of Algo" result values (KNN imputation re-run) in Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D11").Value = -7.135
$ws.Range("B12").Value = 4.951
$ws.Range("D23").Value = -8.321
$ws.Range("B27").Value = 5.583
$ws.Range("D28").Value = -7.949000000000001
$ws.Range("B32").Value = 5.925
$ws.Range("D32").Value = -7.689
$ws.Range("D34").Value = -7.977999999999999
$ws.Range("B36").Value = 8.599
$ws.Range("B38").Value = 6.477000000000001
$ws.Range("D42").Value = -8.395
$ws.Range("B46").Value = 6.189
$ws.Range("D49").Value = -8.330000000000002
$ws.Range("B54").Value = 4.795999999999999
$ws.Range("D54").Value = -8.167999999999997
$ws.Range("B55").Value = 4.763
$ws.Range("B56").Value = 4.654000000000001
$ws.Range("B67").Value = 5.456999999999999
$ws.Range("B69").Value = 5.259
$ws.Range("B72").Value = 5.234999999999999
$ws.Range("D78").Value = -8.208
$ws.Range("D80").Value = -8.058999999999999
$ws.Range("B83").Value = 5.831
$ws.Range("B86").Value = 5.02
$ws.Range("B91").Value = 5.923
$ws.Range("B93").Value = 4.909000000000001
$ws.Range("D97").Value = -7.106
$ws.Range("B99").Value = 4.76
$ws.Range("D99").Value = -8.373999999999999
$ws.Range("D101").Value = -7.776999999999999
$ws.Range("B104").Value = 7.874000000000001
